# Applies the "Cleared up abstract class explanation in inheritance slide"
# commit: a spelling fix + run split on the Practical:Abstraction slide
# (slide 12), plus a resize/reflow + copy edit of the "What is
# inheritance?" slide (slide 14).

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 12 - "Practical: Abstraction and Encapsulation"
#   AddHoby -> AddHobby, and split the following run so that "(string "
#   becomes its own run ahead of "hobby): will add an entry ...".
# ---------------------------------------------------------------------
$s12 = $p.Slides.Item(12)
$content12 = $s12.Shapes.Item(2)
$tr12 = $content12.TextFrame.TextRange

$methodName = $tr12.Characters(374, 7)
if ($methodName.Text -eq "AddHoby") {
    $methodName.Text = "AddHobby"
}

$paramsStart = $tr12.Paragraphs(8, 1).Start + 5 + 8
$firstPart = $tr12.Characters($paramsStart, 8)
if ($firstPart.Text -eq "(string ") {
    # Re-assign the identical text: this forces PowerPoint to split the
    # run in two (same rPr on both halves) instead of leaving one run.
    $firstPart.Text = "(string "
}

# ---------------------------------------------------------------------
# Slide 14 - "What is inheritance?"
# ---------------------------------------------------------------------
$s14 = $p.Slides.Item(14)

# The three decorative "Designer" rectangles got re-generated ids the
# next time the design was touched; only their display Name is
# reachable from automation (Shape.Id is read-only), so update that.
$s14.Shapes.Item(1).Name = "Rectangle 15"
$s14.Shapes.Item(2).Name = "Rectangle 17"
$s14.Shapes.Item(3).Name = "Rectangle 19"

$title = $s14.Shapes.Item(5)
$content14 = $s14.Shapes.Item(6)

# Reflow: title box shrinks, content box grows up into the space and
# down further, and both keep (basically) the same left edge / width.
$title.Left = 38.769330708661414
$title.Top = 40.695787401574805
$title.Width = 242.90114173228346
$title.Height = 99.86492125984252

$content14.Left = 38.7694094488189
$content14.Top = 140.56066929133857
$content14.Width = 242.90114173228346
$content14.Height = 331.0393307086614

$tr14 = $content14.TextFrame.TextRange

# Bump every run in the body from 14pt to 16pt.
$tr14.Font.Size = 16

# Paragraph 2: drop "or partial" before "classes."
$para2 = $tr14.Characters(101, 142)
$para2.Text = "For example, we can make a class for mammals, and have dogs inherit from that class. Often, we will inherit from abstract classes. "

# Paragraph 3: reworded explanation of abstract classes.
$para3start = $tr14.Paragraphs(3, 1).Start
$para3 = $tr14.Characters($para3start, 250)
$rdquo = [char]0x201D
$para3.Text = "It is often useful to make an abstract class, which is an " + $rdquo + "incomplete class"", or to use virtual methods when making a class that will act as a parent class. You can't make a new instance of an abstract class since it's incomplete, but you can inherit from them and build on top of them."
